$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item("TextBox 19")

# Widen the textbox to fit the new subscript "k" (512961 -> 583493 EMU == 40.390629... -> 45.944330... pt)
$shape.Width = 45.94433070866142

$tr = $shape.TextFrame.TextRange
$tr.Text = "acok"

# Format the newly appended "k" as an italic, subscripted Times New Roman run,
# matching the "(zk)" textbox's existing subscript-k formatting.
$kChar = $tr.Characters(4, 1)
$kChar.Font.Italic = $true
$kChar.Font.BaselineOffset = -0.25
$kChar.Font.Name = "Times New Roman"
